$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) - numeric-looking values must stay as text,
# matching the original inlineStr/text storage for this column.
$dValues = @{
    2 = "237.13"
    3 = "21.56"
    4 = "5.462"
    5 = "0.05652"
    6 = "6.498"
    7 = "3.355"
    8 = "1.067"
    9 = "0.7919"
    10 = "0.1399"
    11 = "0.07330"
    13 = "0.02990"
    14 = "0.09254"
    15 = "0.001674"
    16 = "3.262"
    17 = "0.04773"
    18 = "0.0005747"
    19 = "0.006232"
    20 = "0.005098"
    21 = "0.001050"
    22 = "0.0001502"
    23 = "3.899"
    26 = "0.1055"
    27 = "0.0008320"
    41 = "0.006962"
    42 = "0.003504"
    43 = "0.1039"
    44 = "0.009714"
    45 = "0.00005449"
    47 = "0.6761"
    48 = "0.03680"
}

foreach ($row in $dValues.Keys) {
    $cell = $ws.Cells.Item([int]$row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $dValues[$row]
    $cell.Style = "Normal"
}

# Column E (Volume(1h)) - plain text labels.
$eValues = @{
    18 = "17OneONEWorstin24h"
    27 = "26UpBotsUBXTBestin24h"
    42 = "41CEJICEJI"
    48 = "47BOLOBOLO"
}

foreach ($row in $eValues.Keys) {
    $ws.Cells.Item([int]$row, 5).Value = $eValues[$row]
}

# Column G (Hora) - numeric-looking hour values must also stay as text.
$gValues = @{
    2 = "12"
    3 = "12"
    4 = "12"
    5 = "12"
    6 = "12"
    7 = "12"
    8 = "12"
    9 = "12"
    10 = "12"
    11 = "12"
    12 = "12"
    13 = "12"
    14 = "12"
    15 = "12"
    16 = "12"
    17 = "12"
    18 = "12"
    19 = "12"
    20 = "12"
    21 = "12"
    22 = "12"
    23 = "12"
    24 = "12"
    25 = "12"
    26 = "12"
    27 = "12"
    28 = "12"
    29 = "12"
    30 = "12"
    31 = "12"
    32 = "12"
    33 = "12"
    34 = "12"
    35 = "12"
    36 = "12"
    37 = "12"
    38 = "12"
    39 = "12"
    40 = "12"
    41 = "12"
    42 = "12"
    43 = "12"
    44 = "12"
    45 = "12"
    46 = "12"
    47 = "12"
    48 = "12"
    49 = "12"
    50 = "12"
    51 = "12"
}

foreach ($row in $gValues.Keys) {
    $cell = $ws.Cells.Item([int]$row, 7)
    $cell.NumberFormat = "@"
    $cell.Value = $gValues[$row]
    $cell.Style = "Normal"
}
